$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.629.62"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "2.038.73"
$ws.Range("E3").Value = "  +2.71%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'258.04"
$ws.Range("E5").Value = "  +5.02%  "
$ws.Range("D6").Value = "'0.624"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("D7").Value = "'57.88"
$ws.Range("E7").Value = "  -5.72%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.387"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").Value = "'14.82"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "2.339.84"
$ws.Range("E14").Value = "  +2.88%  "
$ws.Range("E15").Value = "  -3.36%  "
$ws.Range("D16").Value = "'21.41"
$ws.Range("E16").Value = "  -4.22%  "
$ws.Range("D17").Value = "'5.37"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").Value = "2.039.80"
$ws.Range("E18").Value = "  +3.00%  "
$ws.Range("D19").Value = "37.547.41"
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("D20").Value = "'70.12"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("D23").Value = "'229.87"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").Value = "'2.65"
$ws.Range("E24").Value = "  +5.50%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("D27").Value = "'9.15"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("D28").Value = "'164.01"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  -5.91%  "
$ws.Range("D30").Value = "'20.08"
$ws.Range("E30").Value = "  +2.44%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").Value = "'0.0665"
$ws.Range("E33").Value = "  +6.67%  "
$ws.Range("E34").Value = "  -2.62%  "
$ws.Range("D35").Value = "'4.52"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").Value = "'2.50"
$ws.Range("E36").Value = "  +9.43%  "
$ws.Range("D37").Value = "'3.48"
$ws.Range("E37").Value = "  +3.46%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("E40").Value = "  -3.21%  "
$ws.Range("E41").Value = "  +4.01%  "
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("D43").Value = "'0.0217"
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("D44").Value = "'1.18"
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D45").Value = "1.412.74"
$ws.Range("E45").Value = "  +2.47%  "
$ws.Range("D46").Value = "'16.14"
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").Value = "'91.21"
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").Value = "'7.41"
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("D50").Value = "'2.87"
$ws.Range("E50").Value = "  +2.15%  "
$ws.Range("E51").Value = "  -1.23%  "
